$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simple Data")

$ws.Range("J2").Value = "1 ave"
$ws.Range("K2").Formula = "=AVERAGE(D2:D66)"

$ws.Range("J3").Value = "2 ave"
$ws.Range("K3").Formula = "=AVERAGE(E2:E66)"

$ws.Range("J4").Value = "3 ave"
$ws.Range("K4").Formula = "=AVERAGE(F2:F66)"

$ws.Range("J5").Value = "4 ave"
$ws.Range("K5").Formula = "=AVERAGE(G2:G66)"

$ws.Range("J6").Value = "5 ave"
$ws.Range("K6").Formula = "=AVERAGE(H2:H66)"

$ws.Range("K7").Select() | Out-Null
